# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
# 1. Insert a new "Player Info" sheet in front of the existing sheets.
# 2. Rename the "MATCH_CARD_LINK" columns to "MATCH_CODE" on the
#    "ODI Batting" and "ODI Bowling" sheets, replacing the full scorecard
#    URL values with just the numeric match code.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "Player Info" sheet, placed before "ODI Batting"
#    NOTE: sheet handles returned by Worksheets.Item(...) track the
#    sheet's *position*, not a fixed identity, so re-resolve handles by
#    name after any operation that reorders/inserts sheets.
# ---------------------------------------------------------------------
$battingSheetForInsert = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheetForInsert)
$playerInfo.Name = "Player Info"

$battingSheet = $wb.Worksheets.Item("ODI Batting")

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Reuse the existing bold/bordered header formatting (Range.Style assignment
# is not wired up in this host, so copy the format from a known header cell).
$battingSheet.Range("A1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)

$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4852"
$playerInfo.Range("B2").Value = "Bodiyabaduge Oshada Piyumal Fernando"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Leg Break"

# ---------------------------------------------------------------------
# 2. "ODI Batting": MATCH_CARD_LINK (col D) -> MATCH_CODE, URL -> code
# ---------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingCodes = @{
    2 = "4261"
    3 = "4264"
    4 = "4269"
    5 = "4271"
    6 = "4272"
    7 = "4375"
    8 = "4450"
    9 = "4471"
}

foreach ($row in $battingCodes.Keys) {
    $cell = $battingSheet.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $battingCodes[$row]
}

# ---------------------------------------------------------------------
# 3. "ODI Bowling": MATCH_CARD_LINK (col B) -> MATCH_CODE, URL -> code
# ---------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingCell = $bowlingSheet.Cells.Item(2, 2)
$bowlingCell.NumberFormat = "@"
$bowlingCell.Value = "4271"

Write-Output "edit complete"
